$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking per correct answer (B11): 5 -> 4
$ws.Range("B11").Value = 4
# Marking penalty per wrong answer (C11): -1 -> -2
$ws.Range("C11").Value = -2

# Recalculated totals based on new marking scheme
# Total marks from right answers (B12): 80 -> 64
$ws.Range("B12").Value = 64
# Total marks lost from wrong answers (C12): -9 -> -18
$ws.Range("C12").Value = -18
# Score summary text (E12): "80 / 140" -> "46 / 112"
$ws.Range("E12").Value = "46 / 112"
